$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 527.4643
$arr[0,1] = 0
$arr[0,2] = 527.4643
$arr[0,3] = 0
$arr[0,4] = 1582.3929
$arr[0,5] = $null
$arr[0,6] = -1918.3929
$ws.Range("H17:N17").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2054
$arr[0,1] = 1496
$arr[0,2] = 3356
$arr[0,3] = 1496
$arr[0,4] = 3356
$arr[0,5] = -1011
$arr[0,6] = -4326
$ws.Range("H28:N28").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 236.26666
$arr[0,1] = 248.59259
$arr[0,2] = 125.333336
$arr[0,3] = 248.59259
$arr[0,4] = 125.333336
$arr[0,5] = -19.59259
$arr[0,6] = -583.333336
$ws.Range("H33:N33").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1032
$arr[0,1] = 0
$arr[0,2] = 1032
$arr[0,3] = 0
$arr[0,4] = 3096
$arr[0,5] = $null
$arr[0,6] = -9230
$ws.Range("H111:N111").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3389.1667
$arr[0,1] = 2505
$arr[0,2] = 3566
$arr[0,3] = 2505
$arr[0,4] = 3566
$arr[0,5] = 749
$arr[0,6] = -10074
$ws.Range("H113:N113").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H114:N114").Value = $arr

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2291.875
$arr[0,1] = 1205.5
$arr[0,2] = 3378.25
$arr[0,3] = 1205.5
$arr[0,4] = 3378.25
$arr[0,5] = -1092.5
$arr[0,6] = -3604.25
$ws.Range("H2:N2").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10842.857
$arr[0,1] = 9990
$arr[0,2] = 12975
$arr[0,3] = 9990
$arr[0,4] = 12975
$arr[0,5] = -9731
$arr[0,6] = -13493
$ws.Range("H23:N23").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 39900
$arr[0,1] = 0
$arr[0,2] = 39900
$arr[0,3] = 0
$arr[0,4] = 39900
$arr[0,5] = $null
$arr[0,6] = -40446
$ws.Range("H37:N37").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 28333.334
$arr[0,1] = 5000
$arr[0,2] = 40000
$arr[0,3] = 5000
$arr[0,4] = 40000
$arr[0,5] = -4685
$arr[0,6] = -40630
$ws.Range("H55:N55").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1046.0555
$arr[0,1] = 786.8461
$arr[0,2] = 1720
$arr[0,3] = 786.8461
$arr[0,4] = 1720
$arr[0,5] = 87.15390000000002
$arr[0,6] = -3468
$ws.Range("H74:N74").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1046.0555
$arr[0,1] = 786.8461
$arr[0,2] = 1720
$arr[0,3] = 3934.2305
$arr[0,4] = 8600
$arr[0,5] = 433.7695000000003
$arr[0,6] = -17336
$ws.Range("H77:N77").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5208
$arr[0,1] = 6629.091
$arr[0,2] = 1300
$arr[0,3] = 6629.091
$arr[0,4] = 1300
$arr[0,5] = -5007.091
$arr[0,6] = -4544
$ws.Range("H102:N102").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2291.875
$arr[0,1] = 1205.5
$arr[0,2] = 3378.25
$arr[0,3] = 1205.5
$arr[0,4] = 3378.25
$arr[0,5] = 1088.5
$arr[0,6] = -7966.25
$ws.Range("H116:N116").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6069.2666
$arr[0,1] = 7830.143
$arr[0,2] = 4528.5
$arr[0,3] = 23490.429
$arr[0,4] = 13585.5
$arr[0,5] = -21040.429
$arr[0,6] = -18485.5
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 24126.6
$arr[0,1] = 0
$arr[0,2] = 24126.6
$arr[0,3] = 0
$arr[0,4] = 24126.6
$arr[0,5] = $null
$arr[0,6] = -33926.6
$ws.Range("H123:N123").Value = $arr

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2291.875
$arr[0,1] = 1205.5
$arr[0,2] = 3378.25
$arr[0,3] = 1205.5
$arr[0,4] = 3378.25
$arr[0,5] = -1091.5
$arr[0,6] = -3606.25
$ws.Range("H3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 14058.667
$arr[0,1] = 15316
$arr[0,2] = 4000
$arr[0,3] = 15316
$arr[0,4] = 4000
$arr[0,5] = -15081
$arr[0,6] = -4470
$ws.Range("H25:N25").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1915.5807
$arr[0,1] = 1776.909
$arr[0,2] = 2254.5557
$arr[0,3] = 1776.909
$arr[0,4] = 2254.5557
$arr[0,5] = 143.0909999999999
$arr[0,6] = -6094.5557
$ws.Range("H107:N107").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H114:N114").Value = $arr

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1398
$arr[0,1] = 1563.3334
$arr[0,2] = 1150
$arr[0,3] = 1563.3334
$arr[0,4] = 1150
$arr[0,5] = -1276.3334
$arr[0,6] = -1724
$ws.Range("H16:N16").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 73986.42999999999
$arr[0,1] = 126976.25
$arr[0,2] = 3333.3333
$arr[0,3] = 126976.25
$arr[0,4] = 3333.3333
$arr[0,5] = -126352.25
$arr[0,6] = -4581.3333
$ws.Range("H62:N62").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 73986.42999999999
$arr[0,1] = 126976.25
$arr[0,2] = 3333.3333
$arr[0,3] = 634881.25
$arr[0,4] = 16666.6665
$arr[0,5] = -631761.25
$arr[0,6] = -22906.6665
$ws.Range("H65:N65").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1999.9166
$arr[0,1] = 1912.375
$arr[0,2] = 2175
$arr[0,3] = 1912.375
$arr[0,4] = 2175
$arr[0,5] = -165.375
$arr[0,6] = -5669
$ws.Range("H105:N105").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 300.33334
$arr[0,1] = 276.42105
$arr[0,2] = 357.125
$arr[0,3] = 276.42105
$arr[0,4] = 357.125
$arr[0,5] = 1643.57895
$arr[0,6] = -4197.125
$ws.Range("H107:N107").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 29371.666
$arr[0,1] = 0
$arr[0,2] = 29371.666
$arr[0,3] = 0
$arr[0,4] = 29371.666
$arr[0,5] = $null
$arr[0,6] = -31451.666
$ws.Range("H109:N109").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1398
$arr[0,1] = 1563.3334
$arr[0,2] = 1150
$arr[0,3] = 1563.3334
$arr[0,4] = 1150
$arr[0,5] = 606.6666
$arr[0,6] = -5490
$ws.Range("H113:N113").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 26429.6
$arr[0,1] = 15000
$arr[0,2] = 27699.555
$arr[0,3] = 15000
$arr[0,4] = 27699.555
$arr[0,5] = -9820
$arr[0,6] = -38059.555
$ws.Range("H141:N141").Value = $arr

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 24.5
$arr[0,1] = 24.5
$arr[0,2] = 0
$arr[0,3] = 73.5
$arr[0,4] = 0
$arr[0,5] = 99.5
$arr[0,6] = $null
$ws.Range("H14:N14").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 270
$arr[0,1] = 270
$arr[0,2] = 0
$arr[0,3] = 2430
$arr[0,4] = 0
$arr[0,5] = 20
$arr[0,6] = $null
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 17545064
$arr[0,1] = 430
$arr[0,2] = 21740520
$arr[0,3] = 1290
$arr[0,4] = 65221560
$arr[0,5] = 3750
$arr[0,6] = -65231640
$ws.Range("H131:N131").Value = $arr

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8133.6665
$arr[0,1] = 0
$arr[0,2] = 8133.6665
$arr[0,3] = 0
$arr[0,4] = 8133.6665
$arr[0,5] = $null
$arr[0,6] = -11877.6665
$ws.Range("H92:N92").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1522.5
$arr[0,1] = 930.25
$arr[0,2] = 1719.9166
$arr[0,3] = 930.25
$arr[0,4] = 1719.9166
$arr[0,5] = 1239.75
$arr[0,6] = -6059.9166
$ws.Range("H113:N113").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H114:N114").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3385.7
$arr[0,1] = 2700
$arr[0,2] = 3754.923
$arr[0,3] = 8100
$arr[0,4] = 11264.769
$arr[0,5] = -5650
$arr[0,6] = -16164.769
$ws.Range("H122:N122").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2466.9412
$arr[0,1] = 1801.3
$arr[0,2] = 3417.8572
$arr[0,3] = 5403.9
$arr[0,4] = 10253.5716
$arr[0,5] = -2873.9
$arr[0,6] = -15313.5716
$ws.Range("H132:N132").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 87226.664
$arr[0,1] = 0
$arr[0,2] = 87226.664
$arr[0,3] = 0
$arr[0,4] = 87226.664
$arr[0,5] = $null
$arr[0,6] = -97586.664
$ws.Range("H140:N140").Value = $arr

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 9890
$arr[0,1] = 9890
$arr[0,2] = 0
$arr[0,3] = 9890
$arr[0,4] = 0
$arr[0,5] = -9754
$arr[0,6] = $null
$ws.Range("H40:N40").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2999.75
$arr[0,1] = 2999.75
$arr[0,2] = 0
$arr[0,3] = 2999.75
$arr[0,4] = 0
$arr[0,5] = -1079.75
$arr[0,6] = $null
$ws.Range("H107:N107").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H108:N108").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 33441340
$arr[0,1] = 0
$arr[0,2] = 33441340
$arr[0,3] = 0
$arr[0,4] = 33441340
$arr[0,5] = $null
$arr[0,6] = -33450518
$ws.Range("H116:N116").Value = $arr

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 28005.2
$arr[0,1] = 0
$arr[0,2] = 28005.2
$arr[0,3] = 0
$arr[0,4] = 28005.2
$arr[0,5] = $null
$arr[0,6] = -28581.2
$ws.Range("H15:N15").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 35000
$arr[0,1] = 0
$arr[0,2] = 35000
$arr[0,3] = 0
$arr[0,4] = 35000
$arr[0,5] = $null
$arr[0,6] = -35138
$ws.Range("H27:N27").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 25333.334
$arr[0,1] = 15000
$arr[0,2] = 30500
$arr[0,3] = 15000
$arr[0,4] = 30500
$arr[0,5] = -14480
$arr[0,6] = -31540
$ws.Range("H54:N54").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 34185.668
$arr[0,1] = 0
$arr[0,2] = 34185.668
$arr[0,3] = 0
$arr[0,4] = 34185.668
$arr[0,5] = $null
$arr[0,6] = -34681.668
$ws.Range("H64:N64").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 34185.668
$arr[0,1] = 0
$arr[0,2] = 34185.668
$arr[0,3] = 0
$arr[0,4] = 34185.668
$arr[0,5] = $null
$arr[0,6] = -35901.668
$ws.Range("H67:N67").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 333999.34
$arr[0,1] = 333999.34
$arr[0,2] = 0
$arr[0,3] = 667998.6800000001
$arr[0,4] = 0
$arr[0,5] = -666937.6800000001
$arr[0,6] = $null
$ws.Range("H81:N81").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 333999.34
$arr[0,1] = 333999.34
$arr[0,2] = 0
$arr[0,3] = 3339993.4
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = -3334689.4
$ws.Range("H84:N84").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 27416.334
$arr[0,1] = 0
$arr[0,2] = 27416.334
$arr[0,3] = 0
$arr[0,4] = 27416.334
$arr[0,5] = $null
$arr[0,6] = -30190.334
$ws.Range("H109:N109").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 18863.637
$arr[0,1] = 0
$arr[0,2] = 18863.637
$arr[0,3] = 0
$arr[0,4] = 18863.637
$arr[0,5] = $null
$arr[0,6] = -21997.637
$ws.Range("H115:N115").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1718.76
$arr[0,1] = 920.8542
$arr[0,2] = 3137.2593
$arr[0,3] = 2762.5626
$arr[0,4] = 9411.777900000001
$arr[0,5] = -232.5626000000002
$arr[0,6] = -14471.7779
$ws.Range("H132:N132").Value = $arr
